# Regenerate the localization-status handoff report: the source markdown
# file was re-handed-off under a new GUID name, producing new xlf targets
# (new content hash) and refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

$newGuid = "c386728a-85df-4b69-8b99-b1a6c9fc7119"
$newHash = "d11b9b670c20e351080879509d2fbc578f30673a"

$newHandoffDate  = "2016-09-06 04:57:44"   # Overview!G2 and de-de!H2
$newZhXlfDate    = "2016-09-06 04:57:31"   # zh-cn!H2

$newMdName     = "$newGuid.md"
$newMdDisplay  = "e2e\$newGuid.md"
$newZhXlfName  = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName  = "$newGuid.$newHash.de-de.xlf"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdDisplay
$wsOverview.Range("G2").Value = $newHandoffDate

foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = $newMdDisplay
}

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("G2").Value = $newZhXlfName
$wsZhCn.Range("H2").Value = $newZhXlfDate

foreach ($h in $wsZhCn.Hyperlinks) {
    $h.TextToDisplay = $newMdName
}

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("G2").Value = $newDeXlfName
$wsDeDe.Range("H2").Value = $newHandoffDate

foreach ($h in $wsDeDe.Hyperlinks) {
    $h.TextToDisplay = $newMdName
}
